$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (A4:G4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Conhecimentos Específicos"
$ws.Range("C4").Value = "Gestão da Qualidade"
$ws.Range("D4").Value = "<b>Trilogia Juran</b><br>`r`nControle da Qualidade"
$ws.Range("E4").Value = "<ul>`r`n`t<li>inspeções contínuas para garantir que os processos estejam sob controle</li>`r`n`t<li>fazer um comparativo entre o desempenho real e o desempenho planejado</li>`r`n</ul>"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Row 5 (A5:G5)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Conhecimentos Específicos"
$ws.Range("C5").Value = "Gestão da Qualidade"
$ws.Range("D5").Value = "<b>Trilogia Juran</b><br>`r`nMelhoria da qualidade"
$ws.Range("E5").Value = "<ul>`r`n`t<li>refinamento proativo dos processos para melhorá-los</li>`r`n`t<li>aprimorar o desempenho rumo a um nível superior de qualidade</li>`r`n</ul>"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
